# fix: inversion sole and prey contam (another time !)
# Corrects the sole/diet (prey) contamination values that had been
# inverted in the uncensored BMF diet table, and recomputes the
# BMF_diet_min/median/max columns accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (B,C,D,E,F,G,H,I,J) target values
$data = @{
  2  = @(0.01, 0.01, 0.13, 0.05, 0.1,  0.2,  0.07, 0.15, 2.82)
  3  = @(0.94, 1.41, 2.79, 0.38, 0.56, 0.92, 1.03, 2.54, 7.33)
  4  = @(0.12, 0.22, 0.37, 0.23, 0.37, 0.81, 0.15, 0.6,  1.66)
  5  = @(0.24, 0.5,  0.92, 2.07, 2.51, 3.93, 0.06, 0.2,  0.45)
  6  = @(8.91, 10.59,13.75,4.63, 6.74, 9.72, 0.92, 1.57, 2.97)
  7  = @(1.1,  1.57, 2.15, 0.46, 1.03, 1.3,  0.84, 1.52, 4.65)
  8  = @(0.66, 1.01, 1.24, 0.55, 0.72, 0.97, 0.68, 1.4,  2.24)
  9  = @(0.89, 1.14, 1.8,  0.42, 0.56, 0.89, 1,    2.03, 4.31)
  10 = @(0.02, 0.02, 0.63, 0.02, 0.48, 0.88, 0.03, 0.05, 26.3)
  11 = @(0.24, 0.41, 0.45, 0.13, 0.18, 0.28, 0.84, 2.28, 3.35)
  12 = @(0.08, 0.08, 0.23, 0.15, 0.24, 0.37, 0.21, 0.35, 1.5)
  13 = @(14.32,17.16,25.8, 10.88,15.45,18.95,0.76, 1.11, 2.37)
}

$columns = @("B", "C", "D", "E", "F", "G", "H", "I", "J")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "$($columns[$i])$row"
        $ws.Range($cellRef).Value = $values[$i]
    }
}
